$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Title banner: MAR_VMSGPRHOST_R005_RELEASE -> APR_VMSGPRHOST_R006_RELEASE
Replace-Text "MAR_VMSGPRHOST_R005_RELEASE" "APR_VMSGPRHOST_R006_RELEASE"

# "RELEASE VMSGPRHOST_R004" -> "RELEASE VMSGPRHOST_R006"  (covers the "IS ALREADY INSTALLED" message)
Replace-Text "RELEASE VMSGPRHOST_R004" "RELEASE VMSGPRHOST_R006"

# "PLEASE EXECUTE THE RELEASE VMSGPRHOST_R004" -> "...R006"
Replace-Text "PLEASE EXECUTE THE RELEASE VMSGPRHOST_R004" "PLEASE EXECUTE THE RELEASE VMSGPRHOST_R006"

# "execute VMS_R004 Components" -> "execute VMS_R006 Components"
Replace-Text "execute VMS_R004 Components" "execute VMS_R006 Components"

# "PLEASE EXECUTE ROLLBACK FOR RELEASE VMSGPRHOST_R004" -> "...R006"
Replace-Text "PLEASE EXECUTE ROLLBACK FOR RELEASE VMSGPRHOST_R004" "PLEASE EXECUTE ROLLBACK FOR RELEASE VMSGPRHOST_R006"

# "ROLLBACK IS ALREADY EXECUTED FOR VMSGPRHOST_R004" -> "...R006"
Replace-Text "ROLLBACK IS ALREADY EXECUTED FOR VMSGPRHOST_R004" "ROLLBACK IS ALREADY EXECUTED FOR VMSGPRHOST_R006"

# "drop the edition R004" -> "drop the edition R000"
Replace-Text "drop the edition R004" "drop the edition R000"

# "set Default R003" -> "set Default R006"
Replace-Text "set Default R003" "set Default R006"

# All remaining FEV_VMSGPRHOST_R004_RELEASE (with or without /ROOTFOLDER suffix) -> APR_VMSGPRHOST_R006_RELEASE
Replace-Text "FEV_VMSGPRHOST_R004_RELEASE" "APR_VMSGPRHOST_R006_RELEASE"

$d.Saved = $false
